$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1729.6666
$ws.Range("I129").Value = 497
$ws.Range("J129").Value = 1841.7273
$ws.Range("K129").Value = 1491
$ws.Range("L129").Value = 5525.1819
$ws.Range("M129").Value = 3509
$ws.Range("N129").Value = -15525.1819

$ws.Range("H135").Value = 1641.7192
$ws.Range("I135").Value = 1558
$ws.Range("J135").Value = 2154.5
$ws.Range("K135").Value = 14022
$ws.Range("L135").Value = 19390.5
$ws.Range("M135").Value = -11487
$ws.Range("N135").Value = -24460.5

$ws.Range("H137").Value = 31251000
$ws.Range("I137").Value = 41667356
$ws.Range("J137").Value = 1933.125
$ws.Range("K137").Value = 125002068
$ws.Range("L137").Value = 5799.375
$ws.Range("M137").Value = -124999518
$ws.Range("N137").Value = -10899.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2149.8367
$ws.Range("I32").Value = 1798.4945
$ws.Range("J32").Value = 6717.2856
$ws.Range("K32").Value = 1798.4945
$ws.Range("L32").Value = 6717.2856
$ws.Range("M32").Value = -1511.4945
$ws.Range("N32").Value = -7291.2856

$ws.Range("H61").Value = 1873.6364
$ws.Range("I61").Value = 1400.289
$ws.Range("J61").Value = 4003.7
$ws.Range("K61").Value = 1400.289
$ws.Range("L61").Value = 4003.7
$ws.Range("M61").Value = -1188.289
$ws.Range("N61").Value = -4427.7

$ws.Range("H63").Value = 22800
$ws.Range("I63").Value = 32775
$ws.Range("K63").Value = 32775
$ws.Range("M63").Value = -32089

$ws.Range("H66").Value = 22800
$ws.Range("I66").Value = 32775
$ws.Range("K66").Value = 163875
$ws.Range("M66").Value = -160443

$ws.Range("H132").Value = 2371.0518
$ws.Range("I132").Value = 1989.4546
$ws.Range("J132").Value = 3570.3572
$ws.Range("K132").Value = 5968.3638
$ws.Range("L132").Value = 10711.0716
$ws.Range("M132").Value = -3438.3638
$ws.Range("N132").Value = -15771.0716

$ws.Range("H136").Value = 1873.6364
$ws.Range("I136").Value = 1400.289
$ws.Range("J136").Value = 4003.7
$ws.Range("K136").Value = 4200.867
$ws.Range("L136").Value = 12011.1
$ws.Range("M136").Value = -1650.867
$ws.Range("N136").Value = -17111.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 16872
$ws.Range("I24").Value = 16872
$ws.Range("K24").Value = 16872
$ws.Range("M24").Value = -16637

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""

$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1711

$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

$ws.Range("H134").Value = 14086342
$ws.Range("I134").Value = 16394908
$ws.Range("J134").Value = 4086.8
$ws.Range("K134").Value = 49184724
$ws.Range("L134").Value = 12260.4
$ws.Range("M134").Value = -49182189
$ws.Range("N134").Value = -17330.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 293.53845
$ws.Range("I22").Value = 241.7
$ws.Range("J22").Value = 466.33334
$ws.Range("K22").Value = 241.7
$ws.Range("L22").Value = 466.33334
$ws.Range("M22").Value = 108.3
$ws.Range("N22").Value = -1166.33334

$ws.Range("H31").Value = 1906.1915
$ws.Range("I31").Value = 1216.3928
$ws.Range("J31").Value = 2922.7368
$ws.Range("K31").Value = 1216.3928
$ws.Range("L31").Value = 2922.7368
$ws.Range("M31").Value = -921.3928000000001
$ws.Range("N31").Value = -3512.7368

$ws.Range("H34").Value = 1906.1915
$ws.Range("I34").Value = 1216.3928
$ws.Range("J34").Value = 2922.7368
$ws.Range("K34").Value = 1216.3928
$ws.Range("L34").Value = 2922.7368
$ws.Range("M34").Value = -1014.3928
$ws.Range("N34").Value = -3326.7368

$ws.Range("H58").Value = 1505.825
$ws.Range("I58").Value = 982.8484999999999
$ws.Range("J58").Value = 3971.2856
$ws.Range("K58").Value = 982.8484999999999
$ws.Range("L58").Value = 3971.2856
$ws.Range("M58").Value = -779.8484999999999
$ws.Range("N58").Value = -4377.2856

$ws.Range("H132").Value = 1652.1384
$ws.Range("I132").Value = 1480.3889
$ws.Range("J132").Value = 2495.2727
$ws.Range("K132").Value = 4441.1667
$ws.Range("L132").Value = 7485.8181
$ws.Range("M132").Value = -1911.1667
$ws.Range("N132").Value = -12545.8181

$ws.Range("H134").Value = 1898.0328
$ws.Range("I134").Value = 1315.0392
$ws.Range("J134").Value = 4871.3
$ws.Range("K134").Value = 3945.1176
$ws.Range("L134").Value = 14613.9
$ws.Range("M134").Value = -1410.1176
$ws.Range("N134").Value = -19683.9

$ws.Range("H136").Value = 1505.825
$ws.Range("I136").Value = 982.8484999999999
$ws.Range("J136").Value = 3971.2856
$ws.Range("K136").Value = 2948.5455
$ws.Range("L136").Value = 11913.8568
$ws.Range("M136").Value = -398.5454999999997
$ws.Range("N136").Value = -17013.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2308.3635
$ws.Range("I136").Value = 1100
$ws.Range("K136").Value = 3300
$ws.Range("M136").Value = 1800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H24").Value = 50000
$ws.Range("I24").Value = 50000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 50000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -49827
$ws.Range("N24").Value = ""

$ws.Range("H102").Value = 2434.4546
$ws.Range("I102").Value = 2447.375
$ws.Range("J102").Value = 2400
$ws.Range("K102").Value = 2447.375
$ws.Range("L102").Value = 2400
$ws.Range("M102").Value = -825.375
$ws.Range("N102").Value = -5644

$ws.Range("H132").Value = 3348.0637
$ws.Range("I132").Value = 3282.8333
$ws.Range("K132").Value = 9848.499899999999
$ws.Range("M132").Value = -7318.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""

$ws.Range("H123").Value = 34564.07
$ws.Range("J123").Value = 34564.07
$ws.Range("L123").Value = 34564.07
$ws.Range("N123").Value = -44364.07

$ws.Range("H132").Value = 2018.0745
$ws.Range("I132").Value = 2153.5073
$ws.Range("J132").Value = 1644.28
$ws.Range("K132").Value = 6460.521900000001
$ws.Range("L132").Value = 4932.84
$ws.Range("M132").Value = -3930.521900000001
$ws.Range("N132").Value = -9992.84

$ws.Range("H136").Value = 14409.284
$ws.Range("I136").Value = 20122.314
$ws.Range("J136").Value = 1741.2609
$ws.Range("K136").Value = 60366.942
$ws.Range("L136").Value = 5223.7827
$ws.Range("M136").Value = -57816.942
$ws.Range("N136").Value = -10323.7827
